$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A13").Value = "Big Beat"
$ws.Range("C10").Value = "Funky"
$ws.Range("F15").Select()
